# Control iD - Vault - LeitorRF.xlsx
# Add two new BOM rows: R/1/5%/04 (0402 1 ohm resistor) and the 7-pin JST
# connector used for the RF reader harness.

$xlCenter = -4108   # xlCenter / xlHAlignCenter

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Control iD Vault")

# Row 29 - Res0402R04025%1 / RC1005J1R0CS / R/1/5%/04, supplied by AVNET
# (same layout/column order used by the surrounding resistor rows).
$ws.Range("A29").Value = "Res0402R04025%1"
$ws.Range("B29").Value = "RC1005J1R0CS"
$ws.Range("C29").Value = "R/1/5%/04"
$ws.Range("E29").Value = "AVNET"
$ws.Range("E29").HorizontalAlignment = $xlCenter

# Row 30 - JST7_HORIZJST7_HORIZ connector / A2001WR-S-7P, supplied by Laien.
$ws.Range("A30").Value = "JST7_HORIZJST7_HORIZ"
$ws.Range("C30").Value = "X/JST/7/H/SMD"
$ws.Range("B30").Value = "A2001WR-S-7P"
$ws.Range("E30").Value = "Laien"
$ws.Range("E30").HorizontalAlignment = $xlCenter

# Leave the selection on C30, matching the author's last-edited cell.
$ws.Range("C30").Select() | Out-Null
